# Extract_Matricula_Model.xlsx edit:
#  - Insert two new leading columns: "Id-Demanda" (A) and "Id-Relacao" (B)
#  - Rename the old "Matrícula" header (now shifted to column F) to "Matricula"
#  - Leave the selection on the new "Matricula" header cell (F1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the current column A; existing content
# (including column widths) shifts two positions to the right.
$ws.Range("A:B").Insert()

# Give the two new header cells the same look (fill/border/alignment) as the
# other header cells by copying the format from the neighbouring header C1.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new header text (B1 first so the shared-string table order
# matches: Id-Relacao then Id-Demanda).
$ws.Range("B1").Value = "Id-Relacao"
$ws.Range("A1").Value = "Id-Demanda"

# The old "Matrícula" header (originally D1) is now at F1; drop the accent.
$ws.Range("F1").Value = "Matricula"

# Match the saved selection shown in the workbook.
$ws.Range("F1").Select()
